# Error Calculations and Plots
# Applies the missing-data re-sampling edit:
#  - row 26 (RM 232) and the original row 28 (SC 92) are removed entirely,
#    shifting all following rows up by two
#  - a handful of column D/E ("D" header in col E) values swap between
#    present and missing (blank) across rows 19, 21, 23 and (after the
#    row shift) rows 26, 27, 29, 33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (header "D") missingness changes in rows 19/21/23 -----------
$ws.Range("E19").Value = -6.5
$ws.Range("E21").ClearContents()
$ws.Range("E23").Value = -7

# --- Remove the two rows that disappear from the data set -----------------
# Row 26 = "RM 232" is dropped outright.
$ws.Rows(26).Delete()
# After the delete above, the row that used to be 28 ("SC 92") is now row 27.
$ws.Rows(27).Delete()

# --- Remaining missingness changes, addressed at their POST-SHIFT rows ----
# Row 26 is now "SC 5": B column ("B" header -> col C) becomes missing.
$ws.Range("C26").ClearContents()
# Row 27 is now "SC 101": col C gets a value, col E becomes missing.
$ws.Range("C27").Value = 10
$ws.Range("E27").ClearContents()
# Row 29 is now "SC 119": col C becomes missing.
$ws.Range("C29").ClearContents()
# Row 33 is now "SC 232": col E gets a value.
$ws.Range("E33").Value = -10.7
